# This edit re-orders the content of a handful of rows in the "Artfynd"
# sheet. The row numbers themselves do not move, but the data that used to
# live in one row now lives in another row of the same small group (the
# groups form simple swaps or short rotations). Columns Y and AA hold a
# plain-text date (e.g. "2026-01-25") that is identical for every row inside
# a given group, so those two columns are intentionally left untouched -
# rewriting them through COM would make Excel "helpfully" reinterpret the
# text as a real date serial number, which is not what the source file has.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowSnapshot($r) {
    # Split the row into three chunks so the date-text columns (Y=25, AA=27)
    # can be skipped while everything else is captured in bulk.
    $p1 = $ws.Range("A$r`:X$r").Value2   # columns A..X
    $p2 = $ws.Range("Z$r").Value2        # column Z
    $p3 = $ws.Range("AB$r`:AY$r").Value2 # columns AB..AY
    return @{ P1 = $p1; P2 = $p2; P3 = $p3 }
}

function Set-RowSnapshot($r, $snap) {
    $ws.Range("A$r`:X$r").Value2 = $snap.P1
    $ws.Range("Z$r").Value2 = $snap.P2
    $ws.Range("AB$r`:AY$r").Value2 = $snap.P3
}

# Each inner array lists the rows of one independent rotation group, in
# cycle order: row cycle[i] ends up with the (pre-edit) content that used to
# sit in row cycle[i+1] (wrapping around at the end of the list).
$cycles = @(
    @(25, 26),
    @(45, 48, 47, 46),
    @(61, 62, 63),
    @(78, 79),
    @(80, 81),
    @(82, 83),
    @(89, 90)
)

foreach ($cycle in $cycles) {
    $n = $cycle.Length

    # Snapshot every row in the group before any of them are overwritten.
    $snapshots = @()
    for ($i = 0; $i -lt $n; $i++) {
        $snapshots += ,(Get-RowSnapshot $cycle[$i])
    }

    # Write the rotated content back.
    for ($i = 0; $i -lt $n; $i++) {
        $srcIndex = ($i + 1) % $n
        Set-RowSnapshot $cycle[$i] $snapshots[$srcIndex]
    }
}
